# Trade #15 closed at 2026-02-17 04:07:38 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 15       # Total Trades
$summary.Range("B9").Value = 33.33    # Win Rate %

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 15        # MarketMaking Trades
$status.Range("G4").Value = 33.33     # MarketMaking Win Rate %

# --- Append new trade row (#15) to "All Trades" and "MarketMaking" sheets ---
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Force text (not date) interpretation for date/time-like strings
    $ws.Range("B16").NumberFormat = "@"
    $ws.Range("C16").NumberFormat = "@"

    $ws.Range("A16").Value = 15
    $ws.Range("B16").Value = "2026-02-17"
    $ws.Range("C16").Value = "04:07:32"
    $ws.Range("D16").Value = "MarketMaking"
    $ws.Range("E16").Value = "UP"
    $ws.Range("F16").Value = 0.19
    $ws.Range("G16").Value = 0.19
    $ws.Range("H16").Value = "CLOSED"
    $ws.Range("I16").Value = 0
    $ws.Range("J16").Value = 0
    $ws.Range("K16").Value = 100.03
    $ws.Range("L16").Value = 0
    $ws.Range("M16").Value = 0
    $ws.Range("N16").Value = 0.6
    $ws.Range("O16").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P16").Value = "early_exit"
    $ws.Range("Q16").Value = 0.11

    # Clear the temporary text-number-format styling so the new cells keep
    # the workbook's default (unstyled) appearance, matching the rest of
    # the sheet.
    $ws.Range("A16:Q16").ClearFormats()
}

Write-Output "edit complete"
